{"js": "// Update the date title and the 25 three-digit-by-one-digit\n// multiplication answers in the table, preserving existing formatting.\n\nconst body = context.document.body;\n\n// --- 1. Title paragraph: \"2025-09-07 Sunday\" -> \"2025-09-08 Monday\" ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\n\nif (titlePara.text.trim() === \"2025-09-07 Sunday\") {\n  titlePara.getRange().insertText(\"2025-09-08 Monday\", \"Replace\");\n}\n\n// --- 2. Table cell answers ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Rows that actually contain the 5 answers each (the rest are spacer rows).\nconst dataRows = [0, 4, 9, 14, 19];\n\n// New values, row-major (5 values per populated row).\nconst newValues = [\n  [\"488\u00d79=4392\", \"694\u00d76=4164\", \"150\u00d74=600\", \"132\u00d76=792\", \"177\u00d73=531\"],\n  [\"484\u00d73=1452\", \"863\u00d77=6041\", \"402\u00d73=1206\", \"604\u00d79=5436\", \"736\u00d72=1472\"],\n  [\"414\u00d75=2070\", \"511\u00d79=4599\", \"617\u00d77=4319\", \"999\u00d75=4995\", \"963\u00d73=2889\"],\n  [\"320\u00d79=2880\", \"800\u00d77=5600\", \"425\u00d73=1275\", \"481\u00d79=4329\", \"584\u00d78=4672\"],\n  [\"847\u00d79=7623\", \"960\u00d74=3840\", \"877\u00d74=3508\", \"439\u00d75=2195\", \"746\u00d72=1492\"],\n];\n\nfor (let r = 0; r < dataRows.length; r++) {\n  const rowIndex = dataRows[r];\n  for (let c = 0; c < 5; c++) {\n    const cell = table.getCell(rowIndex, c);\n    cell.body.paragraphs.load(\"items\");\n    await context.sync();\n    const p = cell.body.paragraphs.items[0];\n    p.getRange().insertText(newValues[r][c], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date title and the 25 three-digit-by-one-digit\n# multiplication answers in the table, preserving existing formatting.\n\n$d = $word.ActiveDocument\n\n# --- 1. Title paragraph: \"2025-09-07 Sunday\" -> \"2025-09-08 Monday\" ---\n$titlePara = $d.Paragraphs.Item(1)\nif ($titlePara.Range.Text.Trim() -eq \"2025-09-07 Sunday\") {\n    $titlePara.Range.Text = \"2025-09-08 Monday\"\n}\n\n# --- 2. Table cell answers ---\n$t = $d.Tables.Item(1)\n\n# Rows that actually contain the 5 answers each (the rest are spacer rows).\n# COM Table.Cell(row, col) is 1-indexed.\n$dataRows = @(1, 5, 10, 15, 20)\n\n# New values, row-major (5 values per populated row).\n$newValues = @(\n    @(\"488\u00d79=4392\", \"694\u00d76=4164\", \"150\u00d74=600\", \"132\u00d76=792\", \"177\u00d73=531\"),\n    @(\"484\u00d73=1452\", \"863\u00d77=6041\", \"402\u00d73=1206\", \"604\u00d79=5436\", \"736\u00d72=1472\"),\n    @(\"414\u00d75=2070\", \"511\u00d79=4599\", \"617\u00d77=4319\", \"999\u00d75=4995\", \"963\u00d73=2889\"),\n    @(\"320\u00d79=2880\", \"800\u00d77=5600\", \"425\u00d73=1275\", \"481\u00d79=4329\", \"584\u00d78=4672\"),\n    @(\"847\u00d79=7623\", \"960\u00d74=3840\", \"877\u00d74=3508\", \"439\u00d75=2195\", \"746\u00d72=1492\")\n)\n\nfor ($r = 0; $r -lt $dataRows.Length; $r++) {\n    $rowIndex = $dataRows[$r]\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $t.Cell($rowIndex, $c)\n        $cell.Range.Text = $newValues[$r][$c - 1]\n    }\n}\n"}
